$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.380786895751953
$ws.Range("B1").Value = 1.147139072418213
$ws.Range("C1").Value = 4.374249935150146
$ws.Range("D1").Value = 2.344269275665283
$ws.Range("E1").Value = 0.7315429449081421
